# Add new front-end design
# Updates the truck log data: renumber IDs, swap "Ojo" company name for
# "Rock", rotate the Material values, move dates forward to 4/14/2022,
# and shift the logged times earlier.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = 3
$ws.Range("B2").Value = "Rock"
$ws.Range("D2").Value = "crushedstone"
$ws.Range("F2").Value = 44665
$ws.Range("G2").Value = 0.41666666666666669

# Row 3
$ws.Range("A3").Value = 4
$ws.Range("B3").Value = "Rock"
$ws.Range("D3").Value = "gravel"
$ws.Range("F3").Value = 44665
$ws.Range("G3").Value = 0.41679398148148145

# Row 4
$ws.Range("A4").Value = 4
$ws.Range("B4").Value = "Rock"
$ws.Range("D4").Value = "gravel"
$ws.Range("F4").Value = 44665
$ws.Range("G4").Value = 0.41689814814814818

# Row 5
$ws.Range("A5").Value = 5
$ws.Range("B5").Value = "Rock"
$ws.Range("D5").Value = "lime"
$ws.Range("F5").Value = 44665
$ws.Range("G5").Value = 0.41724537037037041

# Move the active selection from F6 to E6
$ws.Range("E6").Select() | Out-Null
